# Weekly update: insert 3 new rows (Especial/Primera/Segunda) at the top of the
# Mango / Terminal La Palmera de La Serena data block (row 226) and push the
# existing rows down. The new rows carry a new date/price/origin entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows above the current row 226, shifting rows 226:359 down
# to 229:362 (dimension grows from A1:T359 to A1:T362).
$ws.Rows("226:228").Insert()

# Values shared by all three new rows.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$fecha       = 44438
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100108
$producto    = "Tropicales y subtropicales"
$categoriaId = 100108002
$categoria   = "Mango"
$variedad    = "Sin especificar"
$volumen     = 512
$precioMin   = 8500
$precioMax   = 9000
$precioProm  = 8750
$unidad      = "`$/bandeja 4 kilos"
$origen      = "Brasil"
$precioKg    = 2188
$kgUnidad    = 4

$calidades = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = 226 + $i
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidades[$i]
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad
}
